$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "27.815.93"
$ws.Range('E2').Value = "  +1.34%  "
$ws.Range('D3').Value = "1.763.42"
$ws.Range('E3').Value = "  +1.61%  "
$ws.Range('D4').Value = "'1.005"
$ws.Range('E4').Value = "  +0.31%  "
$ws.Range('D5').Value = "'327.71"
$ws.Range('E5').Value = "  +1.53%  "
$ws.Range('E6').Value = "  +0.21%  "
$ws.Range('D7').Value = "'0.4444"
$ws.Range('E7').Value = "  -1.92%  "
$ws.Range('D8').Value = "'0.3536"
$ws.Range('E8').Value = "  +0.16%  "
$ws.Range('D9').Value = "'41.95"
$ws.Range('E9').Value = "  +1.51%  "
$ws.Range('D10').Value = "'0.07373"
$ws.Range('E10').Value = "  +0.06%  "
$ws.Range('D11').Value = "'1.094"
$ws.Range('E11').Value = "  +1.65%  "
$ws.Range('E12').Value = "  +0.28%  "
$ws.Range('D13').Value = "'20.86"
$ws.Range('E13').Value = "  +2.17%  "
$ws.Range('D14').Value = "'6.002"
$ws.Range('E14').Value = "  +1.41%  "
$ws.Range('D15').Value = "'7.211"
$ws.Range('E15').Value = "  +1.99%  "
$ws.Range('D16').Value = "1.764.33"
$ws.Range('E16').Value = "  +1.81%  "
$ws.Range('D17').Value = "'92.90"
$ws.Range('E17').Value = "  +2.09%  "
$ws.Range('D18').Value = "'0.00001058"
$ws.Range('E18').Value = "  +0.66%  "
$ws.Range('D19').Value = "'0.06411"
$ws.Range('E19').Value = "  +1.21%  "
$ws.Range('E20').Value = "  +0.18%  "
$ws.Range('D21').Value = "'17.05"
$ws.Range('E21').Value = "  +2.70%  "
$ws.Range('D22').Value = "'5.755"
$ws.Range('E22').Value = "  +0.27%  "
$ws.Range('D23').Value = "27.876.33"
$ws.Range('E23').Value = "  +1.40%  "
$ws.Range('D24').Value = "'11.23"
$ws.Range('E24').Value = "  +1.02%  "
$ws.Range('D25').Value = "'2.113"
$ws.Range('E25').Value = "  +2.47%  "
$ws.Range('D26').Value = "'161.38"
$ws.Range('E26').Value = "  -0.05%  "
$ws.Range('D27').Value = "'20.42"
$ws.Range('E27').Value = "  +2.55%  "
$ws.Range('D28').Value = "1.969.94"
$ws.Range('E28').Value = "  +2.28%  "
$ws.Range('D29').Value = "'2.141"
$ws.Range('E29').Value = "  +4.58%  "
$ws.Range('D30').Value = "'124.92"
$ws.Range('E30').Value = "  +0.26%  "
$ws.Range('D31').Value = "'1.101"
$ws.Range('E31').Value = "  +5.77%  "
$ws.Range('D32').Value = "'0.09185"
$ws.Range('E32').Value = "  +0.57%  "
$ws.Range('B33').Value = "HuobiToken"
$ws.Range('C33').Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range('D33').Value = "'3.679"
$ws.Range('E33').Value = "  +0.73%  "
$ws.Range('B34').Value = "Filecoin"
$ws.Range('C34').Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range('D34').Value = "'5.613"
$ws.Range('E34').Value = "  +4.31%  "
$ws.Range('D35').Value = "'11.80"
$ws.Range('E35').Value = "  +1.84%  "
$ws.Range('D36').Value = "'0.06166"
$ws.Range('E36').Value = "  +3.79%  "
$ws.Range('D37').Value = "'0.02280"
$ws.Range('E37').Value = "  +0.63%  "
$ws.Range('D38').Value = "'0.2093"
$ws.Range('E38').Value = "  +2.43%  "
$ws.Range('D39').Value = "'0.6300"
$ws.Range('E39').Value = "  +1.09%  "
$ws.Range('D40').Value = "'4.942"
$ws.Range('E40').Value = "  +1.79%  "
$ws.Range('D41').Value = "'1.188"
$ws.Range('E41').Value = "  -0.15%  "
$ws.Range('D42').Value = "'1.393"
$ws.Range('E42').Value = "  +1.94%  "
$ws.Range('D43').Value = "'7.852"
$ws.Range('E43').Value = "  +1.83%  "
$ws.Range('D44').Value = "'13.22"
$ws.Range('E44').Value = "  +1.76%  "
$ws.Range('D45').Value = "'3.747"
$ws.Range('D46').Value = "'0.5847"
$ws.Range('E46').Value = "  +1.03%  "
$ws.Range('D47').Value = "'122.31"
$ws.Range('E47').Value = "  +0.46%  "
$ws.Range('D48').Value = "'1.949"
$ws.Range('E48').Value = "  +1.39%  "
$ws.Range('D49').Value = "'0.06896"
$ws.Range('E49').Value = "  +0.84%  "
$ws.Range('D50').Value = "'1.131"
$ws.Range('E50').Value = "  +1.78%  "
$ws.Range('D51').Value = "'72.73"
$ws.Range('E51').Value = "  +2.61%  "
